$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet ships protected; temporarily unprotect so the cell values below can be written.
$ws.Unprotect()

# Update the confidential disclosure text with the new "as of" date (2021-07-07 -> 2021-07-08)
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-07-08 for illustrative purposes only and are subject to change."

# Refresh Weight (col D) and Percent Change (col E) figures for each holding row
$ws.Range("D2").Value = 0.0141696930473777
$ws.Range("E2").Value = -0.0001157273463717745
$ws.Range("D3").Value = 0.05612712426660525
$ws.Range("E3").Value = 0.009422222703147076
$ws.Range("D4").Value = 0.01587346704069161
$ws.Range("E4").Value = -0.007862718089990839
$ws.Range("D5").Value = 0.008908479283473544
$ws.Range("E5").Value = 0.0002727024815925017
$ws.Range("D6").Value = 0.01546551602250854
$ws.Range("E6").Value = 0.0004673209159491432
$ws.Range("D7").Value = 0.01955334677743596
$ws.Range("E7").Value = -0.000145985401459936
$ws.Range("D8").Value = 0.004418892389603
$ws.Range("E8").Value = -0.002693861843371059
$ws.Range("D9").Value = 0.006604105674436501
$ws.Range("E9").Value = -0.004469458698890816
$ws.Range("D10").Value = 0.01362554578708127
$ws.Range("E10").Value = -0.02659713168187749
$ws.Range("D11").Value = 0.008242833435731043
$ws.Range("E11").Value = -0.003094606542882405
$ws.Range("D12").Value = 0.01350848078967947
$ws.Range("E12").Value = -0.01402751551119497
$ws.Range("D13").Value = 0.002559335403637681
$ws.Range("E13").Value = -0.01495016611295685
$ws.Range("D14").Value = 0.005538343508726153
$ws.Range("E14").Value = 0.001809408926417566
$ws.Range("D15").Value = 0.01352582037814028
$ws.Range("E15").Value = -0.01725372745621445
$ws.Range("D16").Value = 0.009282753237553213
$ws.Range("E16").Value = -0.01766524363315181
$ws.Range("D17").Value = 0.02205969166992808
$ws.Range("E17").Value = 0.001982283342625291
$ws.Range("D18").Value = 0.008563145132902024
$ws.Range("E18").Value = -0.006702412868632712
$ws.Range("D19").Value = 0.01646234497314702
$ws.Range("E19").Value = -0.001947936957676588
$ws.Range("D20").Value = 0.01432231786629896
$ws.Range("E20").Value = 0.005512679162072764
$ws.Range("D21").Value = 0.006501039889434766
$ws.Range("E21").Value = -0.02592464569650876
$ws.Range("D22").Value = 0.01383231506524179
$ws.Range("E22").Value = -0.01376936316695354
$ws.Range("D23").Value = 0.01846520409211858
$ws.Range("E23").Value = -0.006127614226721589
$ws.Range("D24").Value = 0.008841671762258512
$ws.Range("E24").Value = -0.01384462151394439
$ws.Range("D25").Value = 0.02071919875674847
$ws.Range("E25").Value = -0.001905342580595937
$ws.Range("D26").Value = 0.01366353497300855
$ws.Range("E26").Value = -0.0102212265471856
$ws.Range("D27").Value = 0.02049621347116575
$ws.Range("E27").Value = -0.01232093091478037
$ws.Range("D28").Value = 0.06146231217672751
$ws.Range("E28").Value = -0.00919969564916634
$ws.Range("D29").Value = 0.0186397539248204
$ws.Range("E29").Value = -0.01353013530135294
$ws.Range("D30").Value = 0.03111044060653454
$ws.Range("E30").Value = -0.01339215992503506
$ws.Range("D31").Value = 0.01572676379751442
$ws.Range("E31").Value = -0.01200257197971011
$ws.Range("D32").Value = 0.01351734796976801
$ws.Range("E32").Value = -0.009817312617380836
$ws.Range("D33").Value = 0.01710223955816663
$ws.Range("E33").Value = -0.0141907440552288
$ws.Range("D34").Value = 0.04608771511718161
$ws.Range("E34").Value = -0.01130667172699529
$ws.Range("D35").Value = 0.01026163525862892
$ws.Range("E35").Value = 0
$ws.Range("D36").Value = 0.009926990311451804
$ws.Range("E36").Value = -0.003426124197002278
$ws.Range("D37").Value = 0.00968860892893453
$ws.Range("E37").Value = -0.01104842501175363
$ws.Range("D38").Value = 0.006685914520869818
$ws.Range("E38").Value = -0.01021937593677613
$ws.Range("D39").Value = 0.01122594109325771
$ws.Range("E39").Value = -0.02440251572327046
$ws.Range("D40").Value = 0.01714390315776073
$ws.Range("E40").Value = 0
$ws.Range("D41").Value = 0.01756078209014254
$ws.Range("E41").Value = -0.01057957681692734
$ws.Range("D42").Value = 0.0364122247560493
$ws.Range("E42").Value = -0.009107057969926724
$ws.Range("D43").Value = 0.01121060573043335
$ws.Range("E43").Value = -0.008749356664951136
$ws.Range("D44").Value = 0.02259308899271992
$ws.Range("E44").Value = -0.01412499999999994
$ws.Range("D45").Value = 0.01283269234553433
$ws.Range("E45").Value = 0.002044554244579722
$ws.Range("D46").Value = 0.008022824121547903
$ws.Range("E46").Value = -0.00105603815363664
$ws.Range("D47").Value = 0.01173510550608155
$ws.Range("E47").Value = -0.01201732730914329
$ws.Range("D48").Value = 0.009356818484936604
$ws.Range("E48").Value = -0.02278303540133197
$ws.Range("D49").Value = 0.01641533677185571
$ws.Range("E49").Value = -0.02712722267853585
$ws.Range("D50").Value = 0.008157623479126779
$ws.Range("E50").Value = 0.001265662574357451
$ws.Range("D51").Value = 0.01073056332344825
$ws.Range("E51").Value = -0.04193999354769329
$ws.Range("D52").Value = 0.008417079597881884
$ws.Range("E52").Value = -0.003889197555361434
$ws.Range("D53").Value = 0.009724624256280453
$ws.Range("E53").Value = -0.006513945964850998
$ws.Range("D54").Value = 0.1318425781581191
$ws.Range("E54").Value = 0.0001971608832807004
$ws.Range("D55").Value = 0.04507807126929188
$ws.Range("E55").Value = -0.008932661475034243
$ws.Range("D56").Value = 0.9999999999999998
$ws.Range("E56").Value = -0.007059884804398564

# Restore sheet protection to match the original protected state
# (sheet/objects/scenarios protected; column & row formatting left allowed, matching the source file's sheetProtection flags)
$ws.Protect($null, $true, $true, $true, $false, $false, $true, $true)
